$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing row (row 18) into the new row 19 via copy/paste so
# that shared-string-typed cells (lang_code "eng" in A, and the "TRUE" flag in F)
# keep their original cell typing/styling instead of being re-interpreted by a
# plain value assignment (which would turn "TRUE" into a native boolean).
$ws.Range("A18:F18").Copy()
$ws.Range("A19:F19").PasteSpecial()

# New screen key/name/description value.
$ws.Range("B19").Value = "getFirstIdRoot"
$ws.Range("D19").Value = "getFirstIdRoot"
$ws.Range("E19").Value = "getFirstIdRoot"

# C19 (app_id) must stay a genuine number (10003) even though the column's
# number format is Text ("@"); briefly clear the style so the value is stored
# numerically, then restore the Text format without touching the value again.
$ws.Range("C19").Style = "Normal"
$ws.Range("C19").Value = 10003
$ws.Range("C19").NumberFormat = "@"

$ws.Range("F18").Select()
